# chore: update Sheets via scheduled runner
# Refreshes cached market-price/profit figures (currentAveragePrice*, LevePrice*,
# LeveProfit*) on a handful of leve rows across the per-class sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 8444.333000000001
$ws.Range("I74").Value = 9600
$ws.Range("J74").Value = 6999.75
$ws.Range("K74").Value = 9600
$ws.Range("L74").Value = 6999.75
$ws.Range("M74").Value = -8664
$ws.Range("N74").Value = -8871.75

$ws.Range("H77").Value = 8444.333000000001
$ws.Range("I77").Value = 9600
$ws.Range("J77").Value = 6999.75
$ws.Range("K77").Value = 48000
$ws.Range("L77").Value = 34998.75
$ws.Range("M77").Value = -43320
$ws.Range("N77").Value = -44358.75

$ws.Range("H100").Value = 4004.25
$ws.Range("I100").Value = 10005
$ws.Range("J100").Value = 2004
$ws.Range("K100").Value = 10005
$ws.Range("L100").Value = 2004
$ws.Range("M100").Value = -9464
$ws.Range("N100").Value = -3086

$ws.Range("H113").Value = 8101.95
$ws.Range("I113").Value = 15715.571
$ws.Range("J113").Value = 4002.3076
$ws.Range("K113").Value = 15715.571
$ws.Range("L113").Value = 4002.3076
$ws.Range("M113").Value = -12461.571
$ws.Range("N113").Value = -10510.3076

$ws.Range("H141").Value = 5296.1113
$ws.Range("I141").Value = 5765.952
$ws.Range("J141").Value = 3651.6667
$ws.Range("K141").Value = 17297.856
$ws.Range("L141").Value = 10955.0001
$ws.Range("M141").Value = -12117.856
$ws.Range("N141").Value = -21315.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1389.4706
$ws.Range("I2").Value = 1376.9166
$ws.Range("K2").Value = 1376.9166
$ws.Range("M2").Value = -1263.9166

$ws.Range("H61").Value = 4863.2593
$ws.Range("I61").Value = 5379.4546
$ws.Range("J61").Value = 2592
$ws.Range("K61").Value = 5379.4546
$ws.Range("L61").Value = 2592
$ws.Range("M61").Value = -5167.4546
$ws.Range("N61").Value = -3016

$ws.Range("H74").Value = 5783.136
$ws.Range("I74").Value = 541
$ws.Range("J74").Value = 10151.583
$ws.Range("K74").Value = 541
$ws.Range("L74").Value = 10151.583
$ws.Range("M74").Value = 333
$ws.Range("N74").Value = -11899.583

$ws.Range("H77").Value = 5783.136
$ws.Range("I77").Value = 541
$ws.Range("J77").Value = 10151.583
$ws.Range("K77").Value = 2705
$ws.Range("L77").Value = 50757.915
$ws.Range("M77").Value = 1663
$ws.Range("N77").Value = -59493.915

$ws.Range("H110").Value = 688.5238000000001
$ws.Range("I110").Value = 582.6667
$ws.Range("J110").Value = 953.1667
$ws.Range("K110").Value = 582.6667
$ws.Range("L110").Value = 953.1667
$ws.Range("M110").Value = 1462.3333
$ws.Range("N110").Value = -5043.1667

$ws.Range("H116").Value = 1389.4706
$ws.Range("I116").Value = 1376.9166
$ws.Range("K116").Value = 1376.9166
$ws.Range("M116").Value = 917.0834

$ws.Range("H122").Value = 100000904
$ws.Range("I122").Value = 100000904
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 300002712
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -300000262
$ws.Range("N122").ClearContents()  # HQ profit no longer applicable for this row

$ws.Range("H132").Value = 2405252
$ws.Range("I132").Value = 3206084.5
$ws.Range("J132").Value = 2754.7693
$ws.Range("K132").Value = 9618253.5
$ws.Range("L132").Value = 8264.3079
$ws.Range("M132").Value = -9615723.5
$ws.Range("N132").Value = -13324.3079

$ws.Range("H136").Value = 4863.2593
$ws.Range("I136").Value = 5379.4546
$ws.Range("J136").Value = 2592
$ws.Range("K136").Value = 16138.3638
$ws.Range("L136").Value = 7776
$ws.Range("M136").Value = -13588.3638
$ws.Range("N136").Value = -12876

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1389.4706
$ws.Range("I3").Value = 1376.9166
$ws.Range("K3").Value = 1376.9166
$ws.Range("M3").Value = -1262.9166

$ws.Range("H94").Value = 473
$ws.Range("I94").Value = 439.4375
$ws.Range("J94").Value = 1010
$ws.Range("K94").Value = 439.4375
$ws.Range("L94").Value = 1010
$ws.Range("M94").Value = 11.5625
$ws.Range("N94").Value = -1912

$ws.Range("H99").Value = 600.0625
$ws.Range("I99").Value = 595
$ws.Range("J99").Value = 605.125
$ws.Range("K99").Value = 595
$ws.Range("L99").Value = 605.125
$ws.Range("M99").Value = 903
$ws.Range("N99").Value = -3601.125

$ws.Range("H107").Value = 41668690
$ws.Range("I107").Value = 90911000
$ws.Range("J107").Value = 2114
$ws.Range("K107").Value = 90911000
$ws.Range("L107").Value = 2114
$ws.Range("M107").Value = -90909080
$ws.Range("N107").Value = -5954

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 882.8182
$ws.Range("I16").Value = 858.7143
$ws.Range("J16").Value = 925
$ws.Range("K16").Value = 858.7143
$ws.Range("L16").Value = 925
$ws.Range("M16").Value = -571.7143
$ws.Range("N16").Value = -1499

$ws.Range("H113").Value = 882.8182
$ws.Range("I113").Value = 858.7143
$ws.Range("J113").Value = 925
$ws.Range("K113").Value = 858.7143
$ws.Range("L113").Value = 925
$ws.Range("M113").Value = 1311.2857
$ws.Range("N113").Value = -5265

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 9100553
$ws.Range("I113").Value = 833
$ws.Range("J113").Value = 12512948
$ws.Range("K113").Value = 2499
$ws.Range("L113").Value = 37538844
$ws.Range("M113").Value = -329
$ws.Range("N113").Value = -37543184

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 250001230
$ws.Range("I122").Value = 250001230
$ws.Range("K122").Value = 750003690
$ws.Range("M122").Value = -750001240

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2050.6667
$ws.Range("I40").Value = 1252
$ws.Range("J40").Value = 2450
$ws.Range("K40").Value = 1252
$ws.Range("L40").Value = 2450
$ws.Range("M40").Value = -1116
$ws.Range("N40").Value = -2722

$ws.Range("H61").Value = 3021.4546
$ws.Range("I61").Value = 2356.5
$ws.Range("J61").Value = 3401.4285
$ws.Range("K61").Value = 2356.5
$ws.Range("L61").Value = 3401.4285
$ws.Range("M61").Value = -2154.5
$ws.Range("N61").Value = -3805.4285

$ws.Range("H93").Value = 2289.3333
$ws.Range("I93").Value = 800
$ws.Range("J93").Value = 3034
$ws.Range("K93").Value = 800
$ws.Range("L93").Value = 3034
$ws.Range("M93").Value = 448
$ws.Range("N93").Value = -5530

$ws.Range("H100").Value = 2799.889
$ws.Range("I100").Value = 2449.8333
$ws.Range("J100").Value = 3500
$ws.Range("K100").Value = 2449.8333
$ws.Range("L100").Value = 3500
$ws.Range("M100").Value = -1908.8333
$ws.Range("N100").Value = -4582

$ws.Range("H113").Value = 3021.4546
$ws.Range("I113").Value = 2356.5
$ws.Range("J113").Value = 3401.4285
$ws.Range("K113").Value = 2356.5
$ws.Range("L113").Value = 3401.4285
$ws.Range("M113").Value = -186.5
$ws.Range("N113").Value = -7741.4285

$ws.Range("H122").Value = 54350948
$ws.Range("I122").Value = 58825916
$ws.Range("J122").Value = 41671876
$ws.Range("K122").Value = 176477748
$ws.Range("L122").Value = 125015628
$ws.Range("M122").Value = -176475298
$ws.Range("N122").Value = -125020528

$ws.Range("H132").Value = 5002786
$ws.Range("I132").Value = 8335514.5
$ws.Range("K132").Value = 25006543.5
$ws.Range("M132").Value = -25004013.5

$ws.Range("H136").Value = 2211.2666
$ws.Range("I136").Value = 2346.2444
$ws.Range("J136").Value = 1806.3334
$ws.Range("K136").Value = 7038.733200000001
$ws.Range("L136").Value = 5419.0002
$ws.Range("M136").Value = -4488.733200000001
$ws.Range("N136").Value = -10519.0002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1346.6666
$ws.Range("I122").Value = 1350
$ws.Range("J122").Value = 1340
$ws.Range("K122").Value = 4050
$ws.Range("L122").Value = 4020
$ws.Range("M122").Value = -1600
$ws.Range("N122").Value = -8920

$ws.Range("H136").Value = 15584482
$ws.Range("I136").Value = 12314010
$ws.Range("K136").Value = 36942030
$ws.Range("M136").Value = -36939480
